$d = $word.ActiveDocument

$replacements = @(
    @{old="243×2="; new="857×8="},
    @{old="443×9="; new="330×9="},
    @{old="821×8="; new="835×8="},
    @{old="956×5="; new="881×4="},
    @{old="495×8="; new="934×7="},
    @{old="231×4="; new="330×7="},
    @{old="888×5="; new="603×3="},
    @{old="739×4="; new="987×3="},
    @{old="498×7="; new="647×8="},
    @{old="826×8="; new="612×8="},
    @{old="898×6="; new="482×8="},
    @{old="678×3="; new="151×3="},
    @{old="664×4="; new="336×2="},
    @{old="881×9="; new="581×2="},
    @{old="765×4="; new="969×8="},
    @{old="139×2="; new="946×5="},
    @{old="360×8="; new="187×9="},
    @{old="475×2="; new="449×6="},
    @{old="847×3="; new="398×9="},
    @{old="318×6="; new="691×4="},
    @{old="438×8="; new="202×6="},
    @{old="631×6="; new="880×3="},
    @{old="332×7="; new="765×9="},
    @{old="198×2="; new="344×5="},
    @{old="158×4="; new="394×5="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
